# Update "想去人数" (people interested) counts in column F
# for worksheets "展览" and "全部类型" (sheet1 and sheet4),
# matching the data refresh captured in the commit.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 269
    3  = 279
    4  = 10791
    5  = 9530
    6  = 582
    7  = 1
    8  = 697
    9  = 95
    10 = 8
    12 = 19
    13 = 9505
    15 = 2426
    16 = 30
    17 = 70
    18 = 369
    19 = 10826
    20 = 10732
    22 = 11
    23 = 6
    25 = 12
    26 = 9
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
